# Regenerate instance to have positive average demands during the last periods.
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: bump the per-product AverageDemand (column G) ---
$wsProd = $wb.Worksheets.Item("Productdata")
$wsProd.Range("G2").Value = 49
$wsProd.Range("G3").Value = 21
$wsProd.Range("G4").Value = 35
$wsProd.Range("G5").Value = 70

# Column H (LostSale) holds blank cells in the source file; re-blank them so the
# round-trip save doesn't coerce the empty shared-string cells into a value.
$wsProd.Range("H2:H11").Value = ""

# --- ForecastedAverageDemand sheet: fill in the last periods (rows 9-11, B:E) ---
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsAvg.Range("B9").Value = 70
$wsAvg.Range("C9").Value = 30
$wsAvg.Range("D9").Value = 50
$wsAvg.Range("E9").Value = 100

$wsAvg.Range("B10").Value = 70
$wsAvg.Range("C10").Value = 30
$wsAvg.Range("D10").Value = 50
$wsAvg.Range("E10").Value = 100

$wsAvg.Range("B11").Value = 70
$wsAvg.Range("C11").Value = 30
$wsAvg.Range("D11").Value = 50
$wsAvg.Range("E11").Value = 100

# --- ForcastedStandardDeviation sheet: matching standard deviations for those periods ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")
$wsStd.Range("B9").Value = 7.166424999999998
$wsStd.Range("C9").Value = 3.071324999999999
$wsStd.Range("D9").Value = 5.118874999999999
$wsStd.Range("E9").Value = 10.23775

$wsStd.Range("B10").Value = 8.1997825
$wsStd.Range("C10").Value = 3.5141925
$wsStd.Range("D10").Value = 5.856987499999999
$wsStd.Range("E10").Value = 11.713975

$wsStd.Range("B11").Value = 9.129804249999998
$wsStd.Range("C11").Value = 3.912773249999999
$wsStd.Range("D11").Value = 6.521288749999998
$wsStd.Range("E11").Value = 13.0425775
